$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = 0.2384378511026075
$ws.Range("C7").Value = 2.626518026364868
$ws.Range("D7").Value = 24.30244504454661
$ws.Range("E7").Value = 4.929751012429188
$ws.Range("F7").Value = 4.990078027221052
$ws.Range("G7").Value = 38

# Row 8
$ws.Range("B8").Value = 0.3423512382350242
$ws.Range("C8").Value = 2.428316310422063
$ws.Range("D8").Value = 24.06023025796693
$ws.Range("E8").Value = 4.905122858600683
$ws.Range("F8").Value = 4.960656200552552
$ws.Range("G8").Value = 37

# Row 9
$ws.Range("B9").Value = 0.5109227334801674
$ws.Range("C9").Value = 3.592714756462251
$ws.Range("D9").Value = 42.77397683544757
$ws.Range("E9").Value = 6.540181712723857
$ws.Range("F9").Value = 6.689578288640727
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = -1.467550273398731
$ws.Range("C10").Value = 3.244036734525904
$ws.Range("D10").Value = 27.95265024820572
$ws.Range("E10").Value = 5.28702659802329
$ws.Range("F10").Value = 5.286668009895974
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = -0.8516085912070004
$ws.Range("C11").Value = 2.641658813510031
$ws.Range("D11").Value = 7.990727642251075
$ws.Range("E11").Value = 2.826787512752077
$ws.Range("F11").Value = 3.013612958234995
$ws.Range("G11").Value = 5
